$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Updated AMS [N] measurements (psoas via points / origins adjustment) -
# column D recalculates automatically since it holds =C#/C$5 formulas.
$ws.Range("C2").Value = 114.3553
$ws.Range("C3").Value = 311.6354
$ws.Range("C4").Value = 463.6341
$ws.Range("C5").Value = 562.0191
$ws.Range("C6").Value = 1185.944
$ws.Range("C7").Value = 1154.811
$ws.Range("C8").Value = 2535.934
$ws.Range("C9").Value = 2043.108

# Move the sheet selection from C6 to H9 (matches the saved cursor
# position recorded in the workbook's sheetView).
$ws.Activate()
$ws.Range("H9").Select()
